$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (e.g. "1.000", "29.174.79") that must remain
# plain text, matching the source inline-string cells. Force text format first
# so Excel does not silently coerce numeric-looking values (e.g. "1.000" -> 1).
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D20","D21","D22","D23","D24","D25","D26","D27","D28","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.174.79'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = '1.854.82'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '237.91'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').Value = '0.6905'
$ws.Range('E6').Value = '  -1.40%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '0.07781'
$ws.Range('E8').Value = '  +3.89%  '
$ws.Range('D9').Value = '0.3044'
$ws.Range('E9').Value = '  -1.27%  '
$ws.Range('D10').Value = '23.19'
$ws.Range('E10').Value = '  -2.66%  '
$ws.Range('D11').Value = '0.08076'
$ws.Range('E11').Value = '  -0.58%  '
$ws.Range('D12').Value = '1.856.81'
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('D13').Value = '0.7213'
$ws.Range('E13').Value = '  -0.95%  '
$ws.Range('D14').Value = '5.187'
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').Value = '89.30'
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('D16').Value = '29.190.00'
$ws.Range('E16').Value = '  -1.22%  '
$ws.Range('D17').Value = '5.734'
$ws.Range('E17').Value = '  -3.07%  '
$ws.Range('D18').Value = '0.000007801'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('D20').Value = '234.54'
$ws.Range('E20').Value = '  -3.47%  '
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('D22').Value = '2.111.49'
$ws.Range('E22').Value = '  -1.47%  '
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '7.485'
$ws.Range('E24').Value = '  -2.14%  '
$ws.Range('D25').Value = '161.97'
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').Value = '8.964'
$ws.Range('E26').Value = '  -0.99%  '
$ws.Range('D27').Value = '0.1421'
$ws.Range('E27').Value = '  -4.04%  '
$ws.Range('D28').Value = '18.02'
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('D30').Value = '1.402'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').Value = '4.507'
$ws.Range('E31').Value = '  +2.04%  '
$ws.Range('D32').Value = '1.481'
$ws.Range('E32').Value = '  -1.77%  '
$ws.Range('D33').Value = '4.007'
$ws.Range('E33').Value = '  -1.61%  '
$ws.Range('D34').Value = '0.05193'
$ws.Range('E34').Value = '  -1.50%  '
$ws.Range('D35').Value = '1.179'
$ws.Range('E35').Value = '  -1.99%  '
$ws.Range('D36').Value = '0.7028'
$ws.Range('E36').Value = '  -2.89%  '
$ws.Range('D37').Value = '1.015'
$ws.Range('E37').Value = '  +1.25%  '
$ws.Range('D38').Value = '2.676'
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('D39').Value = '0.01849'
$ws.Range('E39').Value = '  -1.25%  '
$ws.Range('D40').Value = '2.687'
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('D41').Value = '0.9394'
$ws.Range('E41').Value = '  +5.92%  '
$ws.Range('D42').Value = '1.097.82'
$ws.Range('E42').Value = '  +4.32%  '
$ws.Range('D43').Value = '5.960'
$ws.Range('E43').Value = '  +0.58%  '
$ws.Range('D44').Value = '0.4284'
$ws.Range('E44').Value = '  -0.86%  '
$ws.Range('D45').Value = '70.42'
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('D46').Value = '1.000'
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').Value = '102.62'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').Value = '1.793'
$ws.Range('E48').Value = '  +1.94%  '
$ws.Range('D49').Value = '2.007.46'
$ws.Range('E49').Value = '  -1.39%  '
$ws.Range('D50').Value = '9.150'
$ws.Range('E50').Value = '  -1.60%  '
$ws.Range('D51').Value = '6.996'
$ws.Range('E51').Value = '  -3.73%  '
